$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$weekly  = $wb.Worksheets.Item("Weekly")

# --- Summary sheet: updated actuals -----------------------------------
# "Testing" actual hours (F5) increases 39 -> 51; dependent formulas
# (F2 = SUM(F3:F5), G2 = E2-F2, G5 shared formula) recalc automatically.
$summary.Range("F5").Value = 51

# "Research Bluetooth/Location/DB" actual hours (F17) 16 -> 28
$summary.Range("F17").Value = 28

# --- Weekly sheet: new PHASE2 columns (O:S) ----------------------------
$weekly.Columns.Item(15).ColumnWidth = 10.830729166666666

$weekly.Range("O4").Value = "PHASE2: W1"
$weekly.Range("P4").Value = "W2"
$weekly.Range("Q4").Value = "W3"
$weekly.Range("R4").Value = "W4"
$weekly.Range("S4").Value = "W5"

$weekly.Range("O7").Value = 12
$weekly.Range("O9").Value = 12

# --- Selection / active-sheet state ------------------------------------
$summary.Activate() | Out-Null
$summary.Range("F5").Select() | Out-Null

$weekly.Activate() | Out-Null
$weekly.Range("P6").Select() | Out-Null
